$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates: add Profitable (B2), SellPrice (E2), Price Change % (F2),
# and flip Holding (G2) from TRUE to FALSE
$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 77.349997999999999
$ws.Range("F2").Value = -1.0869590792838919
$ws.Range("G2").Value = $false

# New row 3 with remaining Principle balance
$ws.Range("C3").Value = 9891.2999999999993

# Column width adjustments for SellPrice / Price Change % columns
$ws.Columns.Item(5).ColumnWidth = 9.0
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666
